$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "69.120.40"
Set-TextValue $ws.Range("E2") "  -1.40%  "

Set-TextValue $ws.Range("D3") "3.495.58"
Set-TextValue $ws.Range("E3") "  -2.78%  "

Set-TextValue $ws.Range("E4") "  +0.22%  "

Set-TextValue $ws.Range("D5") "573.60"
Set-TextValue $ws.Range("E5") "  -1.26%  "

Set-TextValue $ws.Range("D6") "185.35"
Set-TextValue $ws.Range("E6") "  -3.15%  "

Set-TextValue $ws.Range("D7") "3.485.78"
Set-TextValue $ws.Range("E7") "  -2.97%  "

Set-TextValue $ws.Range("D8") "0.611"
Set-TextValue $ws.Range("E8") "  -3.42%  "

Set-TextValue $ws.Range("D10") "0.187"
Set-TextValue $ws.Range("E10") "  +2.97%  "

Set-TextValue $ws.Range("D11") "0.646"
Set-TextValue $ws.Range("E11") "  -2.97%  "

Set-TextValue $ws.Range("D12") "54.01"
Set-TextValue $ws.Range("E12") "  -3.49%  "

Set-TextValue $ws.Range("D13") "0.0000300"
Set-TextValue $ws.Range("E13") "  -2.44%  "

Set-TextValue $ws.Range("D14") "9.41"
Set-TextValue $ws.Range("E14") "  -2.99%  "

Set-TextValue $ws.Range("D15") "4.061.34"
Set-TextValue $ws.Range("E15") "  -2.86%  "

Set-TextValue $ws.Range("D16") "19.31"
Set-TextValue $ws.Range("E16") "  -3.48%  "

Set-TextValue $ws.Range("D17") "69.159.89"
Set-TextValue $ws.Range("E17") "  -1.37%  "

Set-TextValue $ws.Range("D18") "3.502.97"
Set-TextValue $ws.Range("E18") "  -2.66%  "

Set-TextValue $ws.Range("D19") "12.25"
Set-TextValue $ws.Range("E19") "  -3.58%  "

Set-TextValue $ws.Range("D20") "0.119"
Set-TextValue $ws.Range("E20") "  -1.24%  "

Set-TextValue $ws.Range("D21") "540.27"
Set-TextValue $ws.Range("E21") "  +12.23%  "

Set-TextValue $ws.Range("E22") "  -3.90%  "

Set-TextValue $ws.Range("D23") "18.43"
Set-TextValue $ws.Range("E23") "  -3.95%  "

Set-TextValue $ws.Range("D24") "4.93"
Set-TextValue $ws.Range("E24") "  -2.14%  "

Set-TextValue $ws.Range("D25") "4.42"
Set-TextValue $ws.Range("E25") "  +0.19%  "

Set-TextValue $ws.Range("D26") "93.53"
Set-TextValue $ws.Range("E26") "  -1.96%  "

Set-TextValue $ws.Range("D27") "11.28"
Set-TextValue $ws.Range("E27") "  +1.32%  "

Set-TextValue $ws.Range("E28") "  -2.13%  "

Set-TextValue $ws.Range("D29") "9.08"
Set-TextValue $ws.Range("E29") "  -3.33%  "

Set-TextValue $ws.Range("D30") "31.73"
Set-TextValue $ws.Range("E30") "  -1.49%  "

Set-TextValue $ws.Range("D31") "7.24"
Set-TextValue $ws.Range("E31") "  -6.05%  "

Set-TextValue $ws.Range("E32") "  +2.84%  "

Set-TextValue $ws.Range("D33") "64.34"
Set-TextValue $ws.Range("E33") "  -3.56%  "

Set-TextValue $ws.Range("D34") "0.113"
Set-TextValue $ws.Range("E34") "  -6.03%  "

Set-TextValue $ws.Range("D35") "531.90"
Set-TextValue $ws.Range("E35") "  -9.36%  "

Set-TextValue $ws.Range("D36") "3.07"
Set-TextValue $ws.Range("E36") "  +7.97%  "

Set-TextValue $ws.Range("D37") "37.79"
Set-TextValue $ws.Range("E37") "  -3.21%  "

Set-TextValue $ws.Range("D38") "0.400"
Set-TextValue $ws.Range("E38") "  +0.51%  "

Set-TextValue $ws.Range("E39") "  -0.07%  "

Set-TextValue $ws.Range("D40") "0.0₃0760"
Set-TextValue $ws.Range("E40") "  -5.54%  "

Set-TextValue $ws.Range("D41") "3.36"
Set-TextValue $ws.Range("E41") "  -3.20%  "

Set-TextValue $ws.Range("E42") "  -2.58%  "

Set-TextValue $ws.Range("D43") "3.298.68"
Set-TextValue $ws.Range("E43") "  +2.08%  "

Set-TextValue $ws.Range("D44") "3.04"
Set-TextValue $ws.Range("E44") "  -8.38%  "

Set-TextValue $ws.Range("D45") "2.97"
Set-TextValue $ws.Range("E45") "  -3.14%  "

Set-TextValue $ws.Range("D46") "0.0443"
Set-TextValue $ws.Range("E46") "  -1.58%  "

Set-TextValue $ws.Range("D47") "3.47"
Set-TextValue $ws.Range("E47") "  +4.31%  "

Set-TextValue $ws.Range("D48") "0.133"
Set-TextValue $ws.Range("E48") "  -3.57%  "

Set-TextValue $ws.Range("D49") "8.87"
Set-TextValue $ws.Range("E49") "  -6.08%  "

Set-TextValue $ws.Range("E50") "  -0.01%  "

Set-TextValue $ws.Range("D51") "136.66"
Set-TextValue $ws.Range("E51") "  +2.04%  "
